# feat(matriz-adjacencia-ponderado): cria matriz de adjacencias para grafo ponderado
#
# - Reduce the "vertice" sheet so it only lists vertices A..F (removes G..S).
# - Update the "aresta" sheet so the weight of the C->A edge (row 4, column D)
#   changes from 2 to 1.
# - Restore the active-cell selections that result from this edit.

$wb = $excel.ActiveWorkbook

# --- Sheet "vertice": keep only A1:A6 (A..F), drop rows 7-19 (G..S) ---
$wsVertice = $wb.Worksheets.Item("vertice")
$wsVertice.Range("A7:A19").EntireRow.Delete()
$wsVertice.Range("F4").Select()

# --- Sheet "aresta": weighted adjacency list, CA edge weight 2 -> 1 ---
$wsAresta = $wb.Worksheets.Item("aresta")
$wsAresta.Range("D4").Value = 1
$wsAresta.Range("D4").Select()
